$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; this shifts the existing rows 34-36 down to 35-37
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new weekly data entry
$ws.Cells.Item(34, 1).Value = 8
$ws.Cells.Item(34, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
$ws.Cells.Item(34, 4).Value = 44858
$ws.Cells.Item(34, 5).Value = 4
$ws.Cells.Item(34, 6).Value = 100112026
$ws.Cells.Item(34, 7).Value = "Haba"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 500
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 7500
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(34, 16).Value = 300
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"
